# Fruta / hortaliza, semanal
# Insert a new weekly record for "Vega Monumental Concepción - Acelga" as
# row 210 (pushing the existing rows 210..231 down to 211..232), matching
# the sibling "Primera" quality rows already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 210:231 down by inserting a fresh row at 210.
$ws.Rows.Item(210).Insert()

$row = 210

$ws.Cells.Item($row, 1).Value2 = 11
$ws.Cells.Item($row, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value2 = "Bíobío"
$ws.Cells.Item($row, 4).Value2 = 44637
$ws.Cells.Item($row, 5).Value2 = 8
$ws.Cells.Item($row, 6).Value2 = 100112009
$ws.Cells.Item($row, 7).Value2 = "Acelga"
$ws.Cells.Item($row, 8).Value2 = "Sin especificar"
$ws.Cells.Item($row, 9).Value2 = "Primera"
$ws.Cells.Item($row, 10).Value2 = 350
$ws.Cells.Item($row, 11).Value2 = 600
$ws.Cells.Item($row, 12).Value2 = 650
$ws.Cells.Item($row, 13).Value2 = 621
$ws.Cells.Item($row, 14).Value2 = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item($row, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item($row, 16).Value2 = 621
$ws.Cells.Item($row, 17).Value2 = 1
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"
